$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 260.6316
$ws.Range("I33").Value = 289.84616
$ws.Range("K33").Value = 289.84616
$ws.Range("M33").Value = -60.84616

$ws.Range("H51").Value = 13676.353
$ws.Range("I51").Value = 6999
$ws.Range("K51").Value = 6999
$ws.Range("M51").Value = -6515

$ws.Range("H70").Value = 6474.9
$ws.Range("I70").Value = 2806.5334
$ws.Range("J70").Value = 17480
$ws.Range("K70").Value = 8419.600199999999
$ws.Range("L70").Value = 52440
$ws.Range("M70").Value = -8149.600199999999
$ws.Range("N70").Value = -52980

$ws.Range("H73").Value = 6474.9
$ws.Range("I73").Value = 2806.5334
$ws.Range("J73").Value = 17480
$ws.Range("K73").Value = 8419.600199999999
$ws.Range("L73").Value = 52440
$ws.Range("M73").Value = -7483.600199999999
$ws.Range("N73").Value = -54312

$ws.Range("H116").Value = 4124.4736
$ws.Range("I116").Value = 4150.25
$ws.Range("K116").Value = 4150.25
$ws.Range("M116").Value = -708.25

$ws.Range("H138").Value = 3929.7837
$ws.Range("J138").Value = 4082.6272
$ws.Range("L138").Value = 12247.8816
$ws.Range("N138").Value = -22527.8816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9135.968999999999
$ws.Range("I32").Value = 5119.449
$ws.Range("J32").Value = 21436.562
$ws.Range("K32").Value = 5119.449
$ws.Range("L32").Value = 21436.562
$ws.Range("M32").Value = -4832.449
$ws.Range("N32").Value = -22010.562

$ws.Range("H61").Value = 29414928
$ws.Range("I61").Value = 32261124
$ws.Range("J61").Value = 4233.3335
$ws.Range("K61").Value = 32261124
$ws.Range("L61").Value = 4233.3335
$ws.Range("M61").Value = -32260912
$ws.Range("N61").Value = -4657.3335

$ws.Range("H88").Value = 2359.9473
$ws.Range("J88").Value = 2755.75
$ws.Range("L88").Value = 2755.75
$ws.Range("N88").Value = -3567.75

$ws.Range("H91").Value = 2359.9473
$ws.Range("J91").Value = 2755.75
$ws.Range("L91").Value = 2755.75
$ws.Range("N91").Value = -5563.75

$ws.Range("H97").Value = 1038.5555
$ws.Range("I97").Value = 1137.125
$ws.Range("K97").Value = 1137.125
$ws.Range("M97").Value = -641.125

$ws.Range("H102").Value = 1463.6111
$ws.Range("I102").Value = 1278.3572
$ws.Range("J102").Value = 2112
$ws.Range("K102").Value = 1278.3572
$ws.Range("L102").Value = 2112
$ws.Range("M102").Value = 343.6428000000001
$ws.Range("N102").Value = -5356

$ws.Range("H132").Value = 4767059
$ws.Range("I132").Value = 5886455.5
$ws.Range("K132").Value = 17659366.5
$ws.Range("M132").Value = -17656836.5

$ws.Range("H136").Value = 29414928
$ws.Range("I136").Value = 32261124
$ws.Range("J136").Value = 4233.3335
$ws.Range("K136").Value = 96783372
$ws.Range("L136").Value = 12700.0005
$ws.Range("M136").Value = -96780822
$ws.Range("N136").Value = -17800.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2425.5386
$ws.Range("I20").Value = 2214.2
$ws.Range("J20").Value = 3130
$ws.Range("K20").Value = 2214.2
$ws.Range("L20").Value = 3130
$ws.Range("M20").Value = -1967.2
$ws.Range("N20").Value = -3624

$ws.Range("H64").Value = 104.111115
$ws.Range("I64").Value = 43.666668
$ws.Range("J64").Value = 225
$ws.Range("K64").Value = 43.666668
$ws.Range("L64").Value = 225
$ws.Range("M64").Value = 181.333332
$ws.Range("N64").Value = -675

$ws.Range("H67").Value = 104.111115
$ws.Range("I67").Value = 43.666668
$ws.Range("J67").Value = 225
$ws.Range("K67").Value = 43.666668
$ws.Range("L67").Value = 225
$ws.Range("M67").Value = 736.333332
$ws.Range("N67").Value = -1785

$ws.Range("H134").Value = 10872017
$ws.Range("I134").Value = 11630158
$ws.Range("K134").Value = 34890474
$ws.Range("M134").Value = -34887939

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9975.262000000001
$ws.Range("I31").Value = 4916.1763
$ws.Range("K31").Value = 4916.1763
$ws.Range("M31").Value = -4621.1763

$ws.Range("H34").Value = 9975.262000000001
$ws.Range("I34").Value = 4916.1763
$ws.Range("K34").Value = 4916.1763
$ws.Range("M34").Value = -4714.1763

$ws.Range("H58").Value = 62514976
$ws.Range("I58").Value = 71445304
$ws.Range("J58").Value = 2714
$ws.Range("K58").Value = 71445304
$ws.Range("L58").Value = 2714
$ws.Range("M58").Value = -71445101
$ws.Range("N58").Value = -3120

$ws.Range("H134").Value = 15626813

$ws.Range("H136").Value = 62514976
$ws.Range("I136").Value = 71445304
$ws.Range("J136").Value = 2714
$ws.Range("K136").Value = 214335912
$ws.Range("L136").Value = 8142
$ws.Range("M136").Value = -214333362
$ws.Range("N136").Value = -13242

$ws.Range("H141").Value = 223845.45
$ws.Range("J141").Value = 431320
$ws.Range("L141").Value = 431320
$ws.Range("N141").Value = -441680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 155.33333
$ws.Range("I12").Value = 29.8
$ws.Range("K12").Value = 89.40000000000001
$ws.Range("M12").Value = 83.59999999999999

$ws.Range("H33").Value = 653.4167
$ws.Range("I33").Value = 226.33333
$ws.Range("J33").Value = 1080.5
$ws.Range("K33").Value = 1357.99998
$ws.Range("L33").Value = 6483
$ws.Range("M33").Value = -1074.99998
$ws.Range("N33").Value = -7049

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H64").Value = 6844.143
$ws.Range("J64").Value = 12999
$ws.Range("L64").Value = 38997
$ws.Range("N64").Value = -39537

$ws.Range("H67").Value = 6844.143
$ws.Range("J67").Value = 12999
$ws.Range("L67").Value = 38997
$ws.Range("N67").Value = -40869

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H134").Value = 2545.1
$ws.Range("I134").Value = 2545.1
$ws.Range("K134").Value = 7635.299999999999
$ws.Range("M134").Value = -2565.299999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2616.2144
$ws.Range("I102").Value = 2586.6924
$ws.Range("K102").Value = 2586.6924
$ws.Range("M102").Value = -964.6923999999999

$ws.Range("H113").Value = 106449.8
$ws.Range("I113").Value = 171749.67
$ws.Range("K113").Value = 171749.67
$ws.Range("M113").Value = -169579.67

$ws.Range("H126").Value = 6124.4375
$ws.Range("I126").Value = 5906.7
$ws.Range("K126").Value = 17720.1
$ws.Range("M126").Value = -15250.1

$ws.Range("H132").Value = 4812592.5
$ws.Range("I132").Value = 5212913
$ws.Range("K132").Value = 15638739
$ws.Range("M132").Value = -15636209

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 987.7778
$ws.Range("J46").Value = 1232.1666
$ws.Range("L46").Value = 1232.1666
$ws.Range("N46").Value = -1608.1666

$ws.Range("H55").Value = 567
$ws.Range("J55").Value = 635.3077
$ws.Range("L55").Value = 635.3077
$ws.Range("N55").Value = -981.3077

$ws.Range("H61").Value = 6742.125
$ws.Range("I61").Value = 6017.3076
$ws.Range("J61").Value = 9883
$ws.Range("K61").Value = 6017.3076
$ws.Range("L61").Value = 9883
$ws.Range("M61").Value = -5815.3076
$ws.Range("N61").Value = -10287

$ws.Range("H88").Value = 54499.5
$ws.Range("J88").Value = 54499.5
$ws.Range("L88").Value = 54499.5
$ws.Range("N88").Value = -55355.5

$ws.Range("H91").Value = 54499.5
$ws.Range("J91").Value = 54499.5
$ws.Range("L91").Value = 54499.5
$ws.Range("N91").Value = -57463.5

$ws.Range("H113").Value = 6742.125
$ws.Range("I113").Value = 6017.3076
$ws.Range("J113").Value = 9883
$ws.Range("K113").Value = 6017.3076
$ws.Range("L113").Value = 9883
$ws.Range("M113").Value = -3847.3076
$ws.Range("N113").Value = -14223

$ws.Range("H136").Value = 2215.8948
$ws.Range("I136").Value = 1680.3
$ws.Range("J136").Value = 2811
$ws.Range("K136").Value = 5040.9
$ws.Range("L136").Value = 8433
$ws.Range("M136").Value = -2490.9
$ws.Range("N136").Value = -13533

$ws.Range("H139").Value = 181662.67
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 181662.67
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 181662.67
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -191942.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H96").Value = 4151.76
$ws.Range("I96").Value = 3416.375
$ws.Range("J96").Value = 4497.8237
$ws.Range("K96").Value = 3416.375
$ws.Range("L96").Value = 4497.8237
$ws.Range("M96").Value = -2043.375
$ws.Range("N96").Value = -7243.8237

$ws.Range("H113").Value = 168.5
$ws.Range("I113").Value = 199.16667
$ws.Range("J113").Value = 76.5
$ws.Range("K113").Value = 597.50001
$ws.Range("L113").Value = 229.5
$ws.Range("M113").Value = 1572.49999
$ws.Range("N113").Value = -4569.5

$ws.Range("H132").Value = 14291833
$ws.Range("I132").Value = 17244840
$ws.Range("J132").Value = 18966
$ws.Range("K132").Value = 51734520
$ws.Range("L132").Value = 56898
$ws.Range("M132").Value = -51731990

$ws.Range("H136").Value = 16130757
$ws.Range("J136").Value = 5995.6665
$ws.Range("L136").Value = 17986.9995
$ws.Range("N136").Value = -23086.9995
